$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.699.89'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.63%  '
$ws.Range('E2').ClearFormats()

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.153.34'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E3').ClearFormats()

# Row 4
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E4').ClearFormats()

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.97'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('E5').ClearFormats()

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.70'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E6').ClearFormats()

# Row 7
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E7').ClearFormats()

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.535'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +15.50%  '
$ws.Range('E8').ClearFormats()

# Row 9
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('E9').ClearFormats()

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.433'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.36%  '
$ws.Range('E10').ClearFormats()

# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.99%  '
$ws.Range('E11').ClearFormats()

# Row 12
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.43%  '
$ws.Range('E12').ClearFormats()

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.697.53'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.67%  '
$ws.Range('E13').ClearFormats()

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.05'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.26%  '
$ws.Range('E14').ClearFormats()

# Row 15
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.16%  '
$ws.Range('E15').ClearFormats()

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '58.732.19'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.72%  '
$ws.Range('E16').ClearFormats()

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.77%  '
$ws.Range('E17').ClearFormats()

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.150.99'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.51%  '
$ws.Range('E18').ClearFormats()

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.05'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('E19').ClearFormats()

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.16'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('E20').ClearFormats()

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '372.83'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.47%  '
$ws.Range('E21').ClearFormats()

# Row 22
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('E22').ClearFormats()

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('E23').ClearFormats()

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.07'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('E24').ClearFormats()

# Row 25
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.91%  '
$ws.Range('E25').ClearFormats()

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.167'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('E26').ClearFormats()

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('E27').ClearFormats()

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.26'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +12.91%  '
$ws.Range('E28').ClearFormats()

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0865'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -2.85%  '
$ws.Range('E29').ClearFormats()

# Row 30
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('B30').ClearFormats()
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C30').ClearFormats()
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.88'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('E30').ClearFormats()

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.13'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +3.11%  '
$ws.Range('E31').ClearFormats()

# Row 32
$ws.Range('B32').NumberFormat = '@'
$ws.Range('B32').Value = 'RenderToken'
$ws.Range('B32').ClearFormats()
$ws.Range('C32').NumberFormat = '@'
$ws.Range('C32').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C32').ClearFormats()
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.12'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.48%  '
$ws.Range('E32').ClearFormats()

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.19'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.32%  '
$ws.Range('E33').ClearFormats()

# Row 34
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.95%  '
$ws.Range('E34').ClearFormats()

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '159.20'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.51%  '
$ws.Range('E35').ClearFormats()

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.27'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.73%  '
$ws.Range('E36').ClearFormats()

# Row 37
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.58%  '
$ws.Range('E37').ClearFormats()

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.22'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.90%  '
$ws.Range('E38').ClearFormats()

# Row 39
$ws.Range('B39').NumberFormat = '@'
$ws.Range('B39').Value = 'Stacks'
$ws.Range('B39').ClearFormats()
$ws.Range('C39').NumberFormat = '@'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('C39').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.67'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.41%  '
$ws.Range('E39').ClearFormats()

# Row 40
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'Maker'
$ws.Range('B40').ClearFormats()
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C40').ClearFormats()
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.639.19'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +9.49%  '
$ws.Range('E40').ClearFormats()

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0685'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.19%  '
$ws.Range('E41').ClearFormats()

# Row 42
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +5.32%  '
$ws.Range('E42').ClearFormats()

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.79'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.81%  '
$ws.Range('E43').ClearFormats()

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.710'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('E44').ClearFormats()

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0285'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +5.79%  '
$ws.Range('E45').ClearFormats()

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('E46').ClearFormats()

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.193.92'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('E47').ClearFormats()

# Row 48
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +13.56%  '
$ws.Range('E48').ClearFormats()

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.983'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('E49').ClearFormats()

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.21'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('E50').ClearFormats()

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.36'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.86%  '
$ws.Range('E51').ClearFormats()
